# "made an edit to PM review slides"
#
# Two "Iteration N - X days" header textboxes (on the critical-path-analysis
# slides) get " (2 paths)" appended to their text and are widened to fit the
# new text. The literal point values below are chosen so that, after the
# COM layer's points->EMU (float32) conversion, they land on the exact EMU
# targets from the authored edit (238224/3551455 EMU and 4619100 EMU).

$p = $ppt.ActivePresentation

# --- Slide 9: "Iteration 2 - 6 days" -> "Iteration 2 - 6 days (2 paths)" ---
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(2)
$sh9.Left = 18.757795275590553
$sh9.Width = 279.642136
$sh9.TextFrame.TextRange.Text = "Iteration 2 - 6 days (2 paths)"

# --- Slide 10: "Iteration 3 - 5 days" -> "Iteration 3 - 5 days (2 paths)" ---
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$sh10.Width = 363.708665
$sh10.TextFrame.TextRange.Text = "Iteration 3 - 5 days (2 paths)"
